$d = $word.ActiveDocument

# --- 1. AP-style headline casing: "Software Reuse and Component-based SWE"
#        -> "Software reuse and component-based SWE" -------------------------
# Find the paragraph that holds the homework title and down-case the leading
# letter of "Reuse" and "Component" (only those two words change case; the
# rest of the headline/casing is already AP style).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $par = $d.Paragraphs($i).Range
    if ($par.Text -like "*Software Reuse and Component-based SWE*") {
        for ($w = 1; $w -le $par.Words.Count; $w++) {
            $word1 = $par.Words($w)
            if ($word1.Text -eq "Reuse ") {
                $word1.Characters(1).Text = "r"
            }
            elseif ($word1.Text -eq "Component") {
                $word1.Characters(1).Text = "c"
            }
        }
    }
}

